$d = $word.ActiveDocument
$wdParagraph = 4

function Get-ParagraphIndexByStart($doc, $targetStart) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Start -eq $targetStart) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that currently sits right under
#    the title heading.
# ---------------------------------------------------------------------------
$metaRange = $d.Content.Duplicate
$foundMeta = $metaRange.Find.Execute("Meta description", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundMeta) {
    [void]$metaRange.Expand($wdParagraph)
    $metaIdx = Get-ParagraphIndexByStart $d $metaRange.Start
    if ($metaIdx -gt 0) {
        $d.Paragraphs.Item($metaIdx).Range.Delete()
    }
}

# ---------------------------------------------------------------------------
# 2) Locate the closing "Prompt: Create a feature image ..." paragraph (the
#    very last paragraph in the document) and insert a brand-new bold
#    paragraph right before it, carrying the page title text.
# ---------------------------------------------------------------------------
$promptRange = $d.Content.Duplicate
$foundPrompt = $promptRange.Find.Execute("Prompt: Create a feature image", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundPrompt) {
    [void]$promptRange.Expand($wdParagraph)
    $promptIdx = Get-ParagraphIndexByStart $d $promptRange.Start

    if ($promptIdx -gt 0) {
        [void]$d.Paragraphs.Item($promptIdx).Range.InsertParagraphBefore()

        $titleXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Asgardians for Free - Review of Slot Game | Maximum Winning Potential 7,500x</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        [void]$d.Paragraphs.Item($promptIdx).Range.InsertXML($titleXml)

        # -------------------------------------------------------------
        # 3) Replace the old prompt text with the new meta-description
        #    style copy, while keeping the run's existing (italic)
        #    formatting intact.
        # -------------------------------------------------------------
        $newPromptIdx = $promptIdx + 1
        $oldPromptRange = $d.Paragraphs.Item($newPromptIdx).Range
        $newText = $d.Range($oldPromptRange.Start, $oldPromptRange.End - 1)
        $newText.Text = "Explore the world of Asgard with Endorphina" + [char]0x27 + "s Asgardians Slot Game. Play Asgardians for free and experience the thrill of high-value bonuses and impressive winning potentials of up to 7,500x your bet."
    }
}
